# Updated symbol list (Price column refresh) — GitHub Actions run
# Applies the scraped "Price" (column D) updates for the affected coin rows.
# Values are kept as text (matching the workbook's original inline-string
# cells) by forcing the cell's number format to Text ("@") before writing,
# so Excel doesn't silently re-type them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = [ordered]@{
    2  = "277.13"
    3  = "22.93"
    4  = "6.359"
    5  = "0.06249"
    6  = "3.647"
    7  = "6.645"
    8  = "1.403"
    9  = "0.8328"
    10 = "0.01383"
    12 = "0.08394"
    13 = "0.03521"
    14 = "0.03220"
    15 = "4.079"
    16 = "0.09304"
    17 = "0.001646"
    18 = "0.04728"
    19 = "0.006429"
    20 = "0.005738"
    21 = "0.001079"
    22 = "0.0001502"
    23 = "3.734"
    25 = "0.3327"
    26 = "0.1261"
    28 = "0.0002703"
    41 = "0.007119"
    42 = "0.1175"
    43 = "0.003454"
    44 = "0.01229"
    45 = "0.00005980"
    46 = "0.0009897"
    47 = "0.00000000750"
    48 = "0.7818"
    49 = "0.002480"
    50 = "0.00001399"
    51 = "0.01239"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}
